# Updated cryptos list on Sat Oct 14 12:40:46 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
  "D2"  = "26.899.76";  "E2"  = "  -0.41%  "
  "D3"  = "1.549.70";   "E3"  = "  -0.47%  "
                         "E4"  = "  -0.27%  "
  "D5"  = "206.41";     "E5"  = "  -0.24%  "
                         "E6"  = "  +0.80%  "
                         "E7"  = "  -0.27%  "
  "D8"  = "21.95";      "E8"  = "  +1.12%  "
                         "E9"  = "  -0.44%  "
  "D10" = "0.0596";     "E10" = "  +1.11%  "
  "D11" = "0.0856";     "E11" = "  -0.45%  "
  "D12" = "1.769.32";   "E12" = "  -0.44%  "
  "D13" = "1.550.16";   "E13" = "  -0.31%  "
                         "E14" = "  +0.45%  "
                         "E15" = "  +0.04%  "
  "D16" = "26.892.34";  "E16" = "  -0.32%  "
                         "E17" = "  -0.58%  "
  "D18" = "0.0₃0711";   "E18" = "  +3.12%  "
  "D19" = "216.94";     "E19" = "  +0.91%  "
  "D20" = "7.27";       "E20" = "  -0.34%  "
                         "E21" = "  -0.36%  "
  "D22" = "4.09";       "E22" = "  +1.17%  "
  "D23" = "9.18";       "E23" = "  -0.29%  "
                         "E24" = "  -0.56%  "
  "D25" = "153.16";     "E25" = "  +0.20%  "
  "D26" = "6.66";       "E26" = "  -0.16%  "
  "D27" = "14.94";      "E27" = "  +0.00%  "
                         "E28" = "  +0.53%  "
                         "E29" = "  -0.24%  "
                         "E30" = "  +0.90%  "
                         "E31" = "  -1.43%  "
                         "E32" = "  -0.76%  "
  "D33" = "3.11";       "E33" = "  +3.36%  "
  "D34" = "1.405.89";   "E34" = "  +0.95%  "
                         "E35" = "  +2.25%  "
  "D36" = "0.957";      "E36" = "  -0.94%  "
                         "E37" = "  -0.06%  "
                         "E38" = "  +0.08%  "
                         "E39" = "  +0.64%  "
  "D40" = "0.807";      "E40" = "  -0.62%  "
                         "E41" = "  -0.30%  "
  "D42" = "5.66";       "E42" = "  +3.84%  "
  "D43" = "1.00";       "E43" = "  +1.10%  "
  "D44" = "2.28";       "E44" = "  +0.90%  "
  "D45" = "64.34";      "E45" = "  +0.34%  "
  "D46" = "1.72";       "E46" = "  -1.82%  "
  "D47" = "1.683.74";   "E47" = "  -0.42%  "
  "D48" = "87.05";      "E48" = "  +1.02%  "
                         "E49" = "  +1.41%  "
  "D50" = "0.0₆0101";   "E50" = "  +5.43%  "
                         "E51" = "  -0.19%  "
}

foreach ($addr in $updates.Keys) {
  $rng = $ws.Range($addr)
  # Preserve the cell's existing style while forcing the new value to be
  # stored as text (several of the new prices, e.g. "1.00", would otherwise
  # be auto-coerced into numbers by Excel).
  $origStyle = $rng.Style
  $rng.NumberFormat = "@"
  $rng.Value = $updates[$addr]
  $rng.Style = $origStyle
}
